$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.701.07"
$ws.Range("E2").Value = "  +1.14%  "
$ws.Range("D3").Value = "2.286.96"
$ws.Range("E3").Value = "  +1.35%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'503.59"
$ws.Range("E5").Value = "  +2.10%  "
$ws.Range("D6").Value = "'130.30"
$ws.Range("E6").Value = "  +2.12%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("E8").Value = "  +0.85%  "
$ws.Range("D9").Value = "2.301.03"
$ws.Range("E9").Value = "  +1.43%  "
$ws.Range("D10").Value = "'0.0962"
$ws.Range("E10").Value = "  +1.08%  "
$ws.Range("E11").Value = "  +1.07%  "
$ws.Range("E12").Value = "  +4.30%  "
$ws.Range("E13").Value = "  +5.40%  "
$ws.Range("D14").Value = "'23.22"
$ws.Range("E14").Value = "  +6.66%  "
$ws.Range("D15").Value = "2.693.60"
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("D16").Value = "54.704.48"
$ws.Range("E16").Value = "  +1.34%  "
$ws.Range("E17").Value = "  +1.87%  "
$ws.Range("D18").Value = "2.308.09"
$ws.Range("E18").Value = "  +1.65%  "
$ws.Range("D19").Value = "'10.33"
$ws.Range("E19").Value = "  +3.27%  "
$ws.Range("E20").Value = "  +2.42%  "
$ws.Range("D21").Value = "'306.44"
$ws.Range("E21").Value = "  +2.17%  "
$ws.Range("E22").Value = "  -0.98%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("D24").Value = "'60.91"
$ws.Range("E24").Value = "  -1.81%  "
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("D26").Value = "'0.151"
$ws.Range("E26").Value = "  +2.08%  "
$ws.Range("D28").Value = "'171.93"
$ws.Range("E28").Value = "  +3.00%  "
$ws.Range("D29").Value = "0.0₃0708"
$ws.Range("E29").Value = "  +4.97%  "
$ws.Range("D30").Value = "'1.63"
$ws.Range("E30").Value = "  +1.83%  "
$ws.Range("D31").Value = "'6.05"
$ws.Range("E31").Value = "  +3.33%  "
$ws.Range("E32").Value = "  +3.79%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").Value = "'17.99"
$ws.Range("D35").Value = "'0.996"
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").Value = "'0.938"
$ws.Range("E36").Value = "  +5.39%  "
$ws.Range("E37").Value = "  +1.95%  "
$ws.Range("E38").Value = "  +2.57%  "
$ws.Range("D39").Value = "'36.30"
$ws.Range("E39").Value = "  +1.61%  "
$ws.Range("D40").Value = "'0.377"
$ws.Range("E40").Value = "  +1.80%  "
$ws.Range("E41").Value = "  +2.28%  "
$ws.Range("D42").Value = "'5.09"
$ws.Range("E42").Value = "  +4.99%  "
$ws.Range("E43").Value = "  +2.04%  "
$ws.Range("D44").Value = "'125.79"
$ws.Range("E44").Value = "  +0.61%  "
$ws.Range("D45").Value = "'0.0495"
$ws.Range("E45").Value = "  +2.52%  "
$ws.Range("D46").Value = "'0.0900"
$ws.Range("E46").Value = "  +1.63%  "
$ws.Range("D47").Value = "'246.90"
$ws.Range("E47").Value = "  +5.13%  "
$ws.Range("D48").Value = "'0.551"
$ws.Range("E48").Value = "  +1.87%  "
$ws.Range("D49").Value = "'0.376"
$ws.Range("E49").Value = "  +1.84%  "
$ws.Range("E50").Value = "  +2.72%  "
$ws.Range("E51").Value = "  +0.68%  "
